# Update Week 15 simulation distribution strings on the YDS sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(1.1990050145605002, 0.6443823513378211, 1.498362681117483, 4.164024394887456)"
$ws.Range("C2").Value = "NIG(2.0059339212280443, 1.4612409524903813, 3.367489332121533, 6.661779039357299)"
$ws.Range("D2").Value = "JSU(-1.1406037024274975, 1.1828568557303552, -0.09554257632741492, 2.8120871393497078)"
$ws.Range("E2").Value = "NIG(1.031196866021805, 0.705153327066756, 4.322790061325168, 6.17863226520436)"
